$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Sectors 1A1bc_Other-transformation (row 5) and the 1B* sectors
# (1B1_Fugitive-solid-fuels row 30, 1B2_Fugitive-petr-and-gas row 31,
# 1B2d_Fugitive-other-energy row 32) move from "Energy_Combustion"/"kt"
# to the "process emissions" placeholder activity/units used by the
# other process-emission sectors in this sheet ("GDP"/"B2005USD").
$rows = @(5, 30, 31, 32)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "GDP"
    $ws.Cells.Item($r, 3).Value = "B2005USD"
}
